$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on Price (D) cells whose new value looks like
# a plain number, so Excel does not silently coerce the entered string into
# a numeric cell (which would drop significant trailing zeros and change
# the stored cell type from Text to Number).
$textForceRows = @(
5,6,7,9,10,11,14,16,17,20,22,23,26,27,28,29,30,31,32,35,36,37,38,39,41,42,43,44,45,46,49,50,51
)
foreach ($r in $textForceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "41.523.78"
$ws.Range("E2").Value = "  +4.73%  "

$ws.Range("D3").Value = "2.217.57"
$ws.Range("E3").Value = "  +2.80%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "229.19"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D7").Value = "61.01"
$ws.Range("E7").Value = "  -3.25%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "0.401"
$ws.Range("E9").Value = "  +2.80%  "

$ws.Range("D10").Value = "58.05"
$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("D11").Value = "0.0884"
$ws.Range("E11").Value = "  +4.87%  "

$ws.Range("E12").Value = "  -0.29%  "

$ws.Range("D13").Value = "2.548.66"
$ws.Range("E13").Value = "  +2.91%  "

$ws.Range("D14").Value = "15.55"
$ws.Range("E14").Value = "  -2.09%  "

$ws.Range("E15").Value = "  -1.24%  "

$ws.Range("D16").Value = "0.793"
$ws.Range("E16").Value = "  -1.36%  "

$ws.Range("D17").Value = "5.53"
$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").Value = "2.222.24"

$ws.Range("D19").Value = "41.457.66"
$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").Value = "72.58"
$ws.Range("E20").Value = "  +1.45%  "

$ws.Range("D21").Value = "0.0₃0892"
$ws.Range("E21").Value = "  +5.78%  "

$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").Value = "249.40"
$ws.Range("E23").Value = "  +8.80%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "167.77"
$ws.Range("E28").Value = "  -2.41%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.141"
$ws.Range("E29").Value = "  +1.82%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.43"
$ws.Range("E30").Value = "  -0.94%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "19.84"
$ws.Range("E31").Value = "  +0.24%  "

$ws.Range("D32").Value = "2.56"
$ws.Range("E32").Value = "  -4.47%  "

$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E34").Value = "  +6.24%  "

$ws.Range("D35").Value = "4.61"
$ws.Range("E35").Value = "  +0.85%  "

$ws.Range("D36").Value = "0.0621"
$ws.Range("E36").Value = "  +1.03%  "

$ws.Range("D37").Value = "6.52"
$ws.Range("E37").Value = "  -5.52%  "

$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").Value = "2.36"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "0.000237"
$ws.Range("E41").Value = "  +28.83%  "

$ws.Range("D42").Value = "4.81"
$ws.Range("E42").Value = "  -5.67%  "

$ws.Range("D43").Value = "0.0235"
$ws.Range("E43").Value = "  +4.01%  "

$ws.Range("D44").Value = "8.67"
$ws.Range("E44").Value = "  +10.35%  "

$ws.Range("D45").Value = "0.0981"
$ws.Range("E45").Value = "  +6.87%  "

$ws.Range("D46").Value = "98.70"
$ws.Range("E46").Value = "  -3.91%  "

$ws.Range("E47").Value = "  -0.99%  "

$ws.Range("D48").Value = "1.460.99"
$ws.Range("E48").Value = "  -3.61%  "

$ws.Range("D49").Value = "16.42"
$ws.Range("E49").Value = "  -7.23%  "

$ws.Range("D50").Value = "2.79"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "52.41"
$ws.Range("E51").Value = "  +5.18%  "
